# Update "analyses.xlsx" with the newer source data (adds 2017-2019 and
# refreshes all three data tables on Feuil1): births by weekday, births by
# month, and births by year.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Table 1: Naissances par jour de la semaine (B2:B8) ---
$ws.Range("B2").Value = 5897965
$ws.Range("B3").Value = 6110595
$ws.Range("B4").Value = 5998277
$ws.Range("B5").Value = 6006636
$ws.Range("B6").Value = 6028839
$ws.Range("B7").Value = 5287820
$ws.Range("B8").Value = 4845189

# --- Table 2: Naissances par jour du mois (B11:B22) ---
# C11:C22 hold formulas (=Bn/days) and recompute automatically.
$ws.Range("B11").Value = 3312169
$ws.Range("B12").Value = 3050654
$ws.Range("B13").Value = 3352139
$ws.Range("B14").Value = 3325769
$ws.Range("B15").Value = 3533627
$ws.Range("B16").Value = 3397598
$ws.Range("B17").Value = 3558635
$ws.Range("B18").Value = 3404147
$ws.Range("B19").Value = 3332798
$ws.Range("B20").Value = 3379379
$ws.Range("B21").Value = 3205196
$ws.Range("B22").Value = 3323210

# --- Table 3: Naissances par année (A25:B73 -> extended to A25:B76) ---
# Existing years 1968-2016 (rows 25-73) are unchanged; append 2017-2019.
$ws.Range("A74").Value = 2017
$ws.Range("B74").Value = 730242
$ws.Range("A75").Value = 2018
$ws.Range("B75").Value = 719737
$ws.Range("A76").Value = 2019
$ws.Range("B76").Value = 714029

# Leave the cursor where the author finished editing.
$ws.Activate()
$ws.Range("B23").Select() | Out-Null
